$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 88, shifting existing rows 88-101 down to 89-102.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new weekly price record.
$ws.Range("A88").Value = 7
$ws.Range("B88").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C88").Value = "Ñuble"
$ws.Range("D88").Value = 44522
$ws.Range("E88").Value = 16
$ws.Range("F88").Value = 100112024
$ws.Range("G88").Value = "Choclo"
$ws.Range("H88").Value = "Dulce o Americano"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 100
$ws.Range("K88").Value = 16000
$ws.Range("L88").Value = 17000
$ws.Range("M88").Value = 16500
$ws.Range("N88").Value = "$/malla 60 unidades"
$ws.Range("O88").Value = "Provincia de Limarí"
$ws.Range("P88").Value = 275
$ws.Range("Q88").Value = 60
$ws.Range("R88").Value = "Hortaliza"
